$d = $word.ActiveDocument

# Locate the "Requisitos" Heading2 paragraph that starts the section to remove.
$searchRange = $d.Content
$searchRange.Find.ClearFormatting()
$found = $searchRange.Find.Execute("Requisitos", $false, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)

if ($found) {
    # Delete from the start of that paragraph through to the end of the
    # document body (removes the "Requisitos" heading paragraph and the
    # following bullet-list paragraph with the two requirement lines,
    # leaving the preceding Bibliografia paragraph intact).
    $toRemove = $d.Range($searchRange.Start, $d.Content.End)
    $toRemove.Delete()
}
